# Adds a new weekly price record for Acelga (Femacal de La Calera) as a new
# row inserted right above the current row 177. This pushes the existing
# rows 177-236 down to 178-237 (dimension grows from A1:R236 to A1:R237),
# matching the sole data change described in the commit "Fruta / hortaliza,
# semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 177, shifting rows 177:236 down to 178:237.
$ws.Rows("177:177").Insert()

# Populate the newly inserted row 177 with the new weekly record.
$ws.Cells.Item(177, 1).Value  = 3
$ws.Cells.Item(177, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(177, 3).Value  = "Coquimbo"
$ws.Cells.Item(177, 4).Value  = 44524
$ws.Cells.Item(177, 5).Value  = 5
$ws.Cells.Item(177, 6).Value  = 100112009
$ws.Cells.Item(177, 7).Value  = "Acelga"
$ws.Cells.Item(177, 8).Value  = "Sin especificar"
$ws.Cells.Item(177, 9).Value  = "Primera"
$ws.Cells.Item(177, 10).Value = 290
$ws.Cells.Item(177, 11).Value = 2000
$ws.Cells.Item(177, 12).Value = 2200
$ws.Cells.Item(177, 13).Value = 2110
$ws.Cells.Item(177, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(177, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(177, 16).Value = 352
$ws.Cells.Item(177, 17).Value = 6
$ws.Cells.Item(177, 18).Value = "Hortaliza"
